$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------

# Split the text starting at $start for $len characters into its own run,
# by toggling Bold on then back off. Word (and this engine) materializes a
# run boundary whenever character formatting is applied to a sub-range of
# an existing run, even when the end state matches the original formatting
# (the toggle is a no-op visually/semantically - it only forces the split).
function Split-Run($start, $len) {
    $sub = $d.Range($start, $start + $len)
    $sub.Bold = 1
    $sub.Bold = 0
}

function Assert-StartsWith($range, $expected) {
    $len = $expected.Length
    $actual = $d.Range($range.Start, $range.Start + $len).Text
    if ($actual -ne $expected) {
        throw "Expected paragraph to start with '$expected' but got '$actual'"
    }
}

# ---------------------------------------------------------------------
# 1) "(a) Public Announcement." block: the (i)/(ii) sub-items move from the
#    List2 style to List3, and their leading "(i)"/"(ii)" marker becomes its
#    own run (separate from the rest of the sentence).
# ---------------------------------------------------------------------
$p5 = $d.Paragraphs.Item(5)
Assert-StartsWith $p5.Range "(i)  Report orders"
$p5.Style = "List 3"
Split-Run $p5.Range.Start 3   # "(i)"

$p6 = $d.Paragraphs.Item(6)
Assert-StartsWith $p6.Range "(ii) Submit announcement"
$p6.Style = "List 3"
Split-Run $p6.Range.Start 4   # "(ii)"

# ---------------------------------------------------------------------
# 2) Sample announcement paragraphs: drop the direct widowControl/tabs
#    paragraph formatting in favor of the List1/List2 styles (which already
#    carry equivalent formatting).
# ---------------------------------------------------------------------
$p9 = $d.Paragraphs.Item(9)
Assert-StartsWith $p9.Range "                       (a) Contract award."
$p9.Style = "List 1"

$p10 = $d.Paragraphs.Item(10)
Assert-StartsWith $p10.Range "                       (b) Contract modification."
$p10.Style = "List 1"

$p11 = $d.Paragraphs.Item(11)
Assert-StartsWith $p11.Range "                           (1) This modification adds"
$p11.Style = "List 2"

$p12 = $d.Paragraphs.Item(12)
Assert-StartsWith $p12.Range "                           (2) This modification provides for the purchase"
$p12.Style = "List 2"

$p13 = $d.Paragraphs.Item(13)
Assert-StartsWith $p13.Range "                           (3) This modification provides for the exercise"
$p13.Style = "List 2"

$p14 = $d.Paragraphs.Item(14)
Assert-StartsWith $p14.Range "                           (4) This modification changes"
$p14.Style = "List 2"

# ---------------------------------------------------------------------
# 3) "(D)(5) Miscellaneous data. Include:" paragraph: remove the leading
#    manual line break and split "(D)" from "(5) " into separate runs.
# ---------------------------------------------------------------------
$p15 = $d.Paragraphs.Item(15)
$brk = $d.Range($p15.Range.Start, $p15.Range.Start + 1)
if ([int][char]$brk.Text[0] -ne 11) {
    throw "Expected a manual line break at the start of paragraph 15"
}
$brk.Delete()
$p15 = $d.Paragraphs.Item(15)
Assert-StartsWith $p15.Range "(D)(5) "
Split-Run $p15.Range.Start 3   # "(D)"

# ---------------------------------------------------------------------
# 4) (a)-(d) "Include:" sub-items: drop direct widowControl/tabs formatting
#    (or add a pPr, for the (c) item which previously had none at all) in
#    favor of the List1 style.
# ---------------------------------------------------------------------
$p16 = $d.Paragraphs.Item(16)
Assert-StartsWith $p16.Range "                  (a) a statement that the information"
$p16.Style = "List 1"

$p17 = $d.Paragraphs.Item(17)
Assert-StartsWith $p17.Range "                  (b) any areas of sensitivity"
$p17.Style = "List 1"

$p18 = $d.Paragraphs.Item(18)
Assert-StartsWith $p18.Range "                  (c) indication of appropriate coordination"
$p18.Style = "List 1"

$p19 = $d.Paragraphs.Item(19)
Assert-StartsWith $p19.Range "                  (d) the estimated period of performance"
$p19.Style = "List 1"

# ---------------------------------------------------------------------
# 5) 5205.404-1 Release procedures -> "(a) Application." paragraph.
# ---------------------------------------------------------------------
$p26 = $d.Paragraphs.Item(26)
Assert-StartsWith $p26.Range "  (a) "
$p26.Style = "List 1"
